# Update countries & provincias Spain
# Applies the "paises.xlsx" data refresh: updated timestamp, refreshed
# case counts for several countries, and two countries that moved row
# position in the source feed (Eslovaquia/Mauritania and
# Montserrat/Islas Malvinas swapped alphabetical-ish ordering).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "last refreshed" timestamp banner in A1.
$ws.Range("A1").Value = "Datos actualizados a 24 de Septiembre de 2020 a las 11:14"

# 2. Refreshed per-country counters (B=Casos totales, C=Nuevos casos,
#    D=Casos activos, E=Recuperados, F=Casos criticos, G=Muertes hoy,
#    H=Muertes). Only cells whose values actually changed are touched.

# Row 24: Filipinas
$ws.Range("B24").Value = 296755
$ws.Range("C24").Value = 2180
$ws.Range("D24").Value = 231928
$ws.Range("E24").Value = 59700
$ws.Range("G24").Value = 36
$ws.Range("H24").Value = 5127

# Row 26: Indonesia
$ws.Range("B26").Value = 262022
$ws.Range("C26").Value = 4634
$ws.Range("D26").Value = 191853
$ws.Range("E26").Value = 60064
$ws.Range("G26").Value = 128
$ws.Range("H26").Value = 10105

# Row 29: Canada
$ws.Range("B29").Value = 147756
$ws.Range("C29").Value = 3
$ws.Range("E29").Value = 10725

# Row 42: Oman
$ws.Range("B42").Value = 95907
$ws.Range("C42").Value = 568
$ws.Range("D42").Value = 86765
$ws.Range("E42").Value = 8257
$ws.Range("G42").Value = 10
$ws.Range("H42").Value = 885

# Row 67: Austria
$ws.Range("B67").Value = 40816
$ws.Range("C67").Value = 832
$ws.Range("D67").Value = 31661
$ws.Range("E67").Value = 8372
$ws.Range("G67").Value = 6
$ws.Range("H67").Value = 783

# Row 77: El Salvador
$ws.Range("B77").Value = 28201
$ws.Range("C77").Value = 247
$ws.Range("D77").Value = 22326
$ws.Range("E77").Value = 5052
$ws.Range("G77").Value = 4
$ws.Range("H77").Value = 823

# Row 89: Croacia
$ws.Range("B89").Value = 15572
$ws.Range("C89").Value = 232
$ws.Range("D89").Value = 14111
$ws.Range("E89").Value = 1200
$ws.Range("G89").Value = 4
$ws.Range("H89").Value = 261

# Row 104: Finlandia
$ws.Range("B104").Value = 9379
$ws.Range("C104").Value = 91
$ws.Range("E104").Value = 1186

# Rows 110/111: the feed re-ordered these two countries. Row 110 used to
# hold Mauritania and row 111 Eslovaquia; now row 110 holds the refreshed
# Eslovaquia data and row 111 holds the (unchanged) old Mauritania data.
$ws.Range("A110").Value = "Eslovaquia"
$ws.Range("B110").Value = 7629
$ws.Range("C110").Value = 360
$ws.Range("D110").Value = 3978
$ws.Range("E110").Value = 3610
$ws.Range("F110").Value = 0
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 41

$ws.Range("A111").Value = "Mauritania"
$ws.Range("B111").Value = 7425
$ws.Range("C111").Value = 0
$ws.Range("D111").Value = 7028
$ws.Range("E111").Value = 236
$ws.Range("F111").Value = 0
$ws.Range("G111").Value = 0
$ws.Range("H111").Value = 161

# Row 140: Reunion
$ws.Range("E140").Value = 1008
$ws.Range("H140").Value = 11

# Row 161: Letonia
$ws.Range("B161").Value = 1594
$ws.Range("C161").Value = 22
$ws.Range("E161").Value = 310

# Rows 214/215: same re-ordering pattern as 110/111, this time for
# Islas Malvinas / Montserrat.
$ws.Range("A214").Value = "Montserrat"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 12
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 1

$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("B215").Value = 13
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 13
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 0
